$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

# Row 2 (Malgorzata Zdunik): clear Student ID cell text entirely
$tbl.Cell(2, 2).Range.Delete()

# Row 2: Email -> "@gold.ac.uk" (remove "mzdun001" prefix)
$d.Content.Find.Execute("mzdun001@gold.ac.uk", $false, $false, $false, $false, $false, $true, 1, $false, "@gold.ac.uk", 2)

# Row 2: Work done cell -> single merged run text (no proofErr split)
$d.Content.Find.Execute("Part1, Part2: kNN, Decision Tree, SVM", $false, $false, $false, $false, $false, $true, 1, $false, "Part1, Part2: kNN, Decision Tree, SVM", 2)

# Row 3 (Katherine Knox): clear Student ID cell text entirely
$tbl.Cell(3, 2).Range.Delete()

# Row 3: Email -> "@gold.ac.uk" (remove "kknox001" prefix)
$d.Content.Find.Execute("kknox001@gold.ac.uk", $false, $false, $false, $false, $false, $true, 1, $false, "@gold.ac.uk", 2)

# Row 3: add a _GoBack bookmark at the start of the email cell's paragraph
$emailRange = $tbl.Cell(3, 3).Range
$emailRange.Collapse(1)
$d.Bookmarks.Add("_GoBack", $emailRange)
